# SI_4.pptx edit: refresh the "Updated" date footer field, and tidy up the
# "R = AA.mm(P1, P2)" code sample so it is one run instead of three.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$newDate = "04/12/2025"

# ---------------------------------------------------------------------
# 1) Update the datetimeFigureOut footer placeholder text on the slide
#    master and on every slide layout (they each carry their own copy
#    of the placeholder).
# ---------------------------------------------------------------------
$mp = $m.Shapes.Placeholders
for ($j = 1; $j -le $mp.Count; $j++) {
    $ph = $mp.Item($j)
    if ($ph.PlaceholderFormat.Type -eq 16) {
        $ph.TextFrame.TextRange.Text = $newDate
    }
}

for ($i = 1; $i -le $m.CustomLayouts.Count; $i++) {
    $lay = $m.CustomLayouts.Item($i)
    $placeholders = $lay.Shapes.Placeholders
    for ($j = 1; $j -le $placeholders.Count; $j++) {
        $ph = $placeholders.Item($j)
        if ($ph.PlaceholderFormat.Type -eq 16) {
            $ph.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------
# 2) Collapse the "R = ", "AA.mm", "(P1, P2" runs on slide 1 into a
#    single run reading "R = AA.mm(P1, P2)".
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$group = $slide.Shapes.Item(1)

for ($i = 1; $i -le $group.GroupItems.Count; $i++) {
    $shape = $group.GroupItems.Item($i)
    if ($shape.Name -eq "Rounded Rectangle 107") {
        $tr = $shape.TextFrame.TextRange
        $fullText = $tr.Text
        $start = $fullText.IndexOf("R = AA.mm(P1, P2")
        if ($start -ge 0) {
            $chars = $tr.Characters($start + 1, 16)
            $chars.Text = "R = AA.mm(P1, P2)"
        }
    }
}
